$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Puzzle")
$ws.Activate()
$ws.Range("A1:I1").Insert(-4121)
